# Fix wrong file name when export data to pdf
# -> populate the previously-empty "exported school" worksheet with the
#    actual school rows that should have been written out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Truong trung hoc Chuyen Nguyen Binh Khiem (Vinh Long) ---
$ws.Range("A2").Value = "nbk-vl"
$ws.Range("B2").Value = "Truong trung hoc Chuyen Nguyen Binh Khiem"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Vinh Long"
# Teacher's CMND is a numeric-looking id that must stay text, like the
# header column already is - force text entry then drop the temporary
# number-format override so the cell keeps the sheet's default style.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "337829999"
$ws.Range("E2").ClearFormats()

# --- Row 3: Truong trung hoc Chuyen Nguyen Binh Khiem (Quang Ngai) ---
$ws.Range("A3").Value = "nbk-qn"
$ws.Range("B3").Value = "Truong trung hoc Chuyen Nguyen Binh Khiem"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Quang Ngai"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "334442222"
$ws.Range("E3").ClearFormats()

# Columns A, B and D now hold much longer values than before, so widen
# them accordingly (mirrors the auto-sized column widths written by the
# exporter).
$ws.Columns.Item(1).ColumnWidth = 6.5
$ws.Columns.Item(2).ColumnWidth = 40.83333333333333
$ws.Columns.Item(4).ColumnWidth = 10.666666666666666
